$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 55362.4108122222
$ws.Range("C2").Value = 33697.043609876
$ws.Range("D2").Value = 2193.1172996027399

$ws.Range("B3").Value = 69820.805069573005
$ws.Range("C3").Value = 46807.970267315301
$ws.Range("D3").Value = 845.64969969155095

$ws.Range("B4").Value = 80646.822443359095
$ws.Range("C4").Value = 57222.045606976899
$ws.Range("D4").Value = 433.707665566819

$ws.Range("B5").Value = 92969.225991601095
$ws.Range("C5").Value = 69300.607644599106
$ws.Range("D5").Value = 189.86615494682499

$ws.Range("B6").Value = 100453.238144067
$ws.Range("C6").Value = 76723.548080766705
$ws.Range("D6").Value = 128.79443864645

$ws.Range("B11").Select()
